$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in Horas (D column) estimates for existing rows 10 and 11 ---
$ws.Range("D10").Value = 6
$ws.Range("D11").Value = 4

# --- Insert a new test case "Crear Rol" before the current row 12, pushing
#     "Asignar roles a usuario" (old row 12) down to row 13 and
#     "Asignar Permisos a Rol" (old row 13) down to row 14 (which was blank). ---

# First, move old row13 ("Asignar Permisos a Rol") content into row14,
# copying the cell formatting (s=5) from row13 onto row14 so no new style
# is introduced (row14 cells were style 1/unformatted before).
[void]$ws.Range("B13:F13").Copy()
[void]$ws.Range("B14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Asignar Permisos a Rol"
$ws.Range("C14").Value = "Asignar permisos a un rol"
$ws.Range("D14").Value = 4

# Now move old row12 ("Asignar roles a usuario") content into row13.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Asignar roles a usuario"
$ws.Range("C13").Value = "Asignar roles a un usuario para que pueda utilizar funcionalidad del sistema"
$ws.Range("D13").Value = 5
$ws.Rows.Item(13).RowHeight = 26.25

# Finally, write the new "Crear Rol" case into row12.
$ws.Range("B12").Value = "Crear Rol"
$ws.Range("C12").Value = "Crear un nuevo rol en el sistema"
$ws.Range("D12").Value = 4
$ws.Rows.Item(12).RowHeight = 15.75

# --- Row 21 switches from the "s=2" blank style to the "s=1" blank style
#     used by rows 14-20, matched by copying the format from row 20. ---
[void]$ws.Range("A20:F20").Copy()
[void]$ws.Range("A21:F21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Extend the trailing blank rows by one (new row 30), matching the
#     "s=2" blank style used by rows 22-29. ---
[void]$ws.Range("A29:F29").Copy()
[void]$ws.Range("A30:F30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the saved view: no more frozen/scrolled topLeftCell, and the
#     active selection moves from C10 to B7. ---
[void]$ws.Range("B7").Select()
